# Fix the typo "Alfa Romea" -> "Alfa Romeo" in the TestData_Cars workbook
# (commit message: "Add example for Excel Import")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = "Alfa Romeo"
